$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns; rows 39-42 also get Coin (B) and Link (C) swapped

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.917.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.80%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.665.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.10%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.64%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.91%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.49%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3633'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.72'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.19%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3272'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.30%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.135'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.59%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07098'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.57%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9998'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.60%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.054'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.47%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.63%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.655.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.77%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.614'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.59%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001049'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.08%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06603'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.18%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9994'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.43%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '79.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.29%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.927'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.28%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.73%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.67%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.879.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.89%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.452'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.435'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.22%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.40%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.72%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.838.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.52%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.43%  '

# Row 31
$ws.Range("E31").Value = '  +7.66%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.089'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.746'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.83%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08467'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.19%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.650'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.34%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.02%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.281'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.77%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.185'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.63%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06174'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.41%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02274'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.24%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.314'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.01%  '

# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2077'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.70%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.43%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5941'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.06%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.52'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.88%  '

# Row 46
$ws.Range("E46").Value = '  +2.05%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5636'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.89%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.64%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.951'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.37%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06992'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.78%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.191'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.10%  '
